# Subo avances en el product backlog
# Adds three new bug rows (40, 41, 42 -> ids 42/43 placeholders become rows 44/45)
# to the "LISTA DE BUGS" sheet, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 42 : id 40 "Empate de Partidos"
# ---------------------------------------------------------------------------
$ws.Cells.Item(24,2).Copy() | Out-Null
$ws.Cells.Item(42,2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(42,2).Value = "Empate de Partidos "

$ws.Cells.Item(24,3).Copy() | Out-Null
$ws.Cells.Item(42,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(42,3).Value = "No permite empates cuando hay dos fases: TCT y eliminatorio"

$ws.Cells.Item(42,4).Value = "Pau"
$ws.Cells.Item(42,5).Value = "admin/fechas.aspx"

$ws.Cells.Item(24,6).Copy() | Out-Null
$ws.Cells.Item(42,6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(42,6).Value = "CORREGIDO"

$ws.Rows.Item(42).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 43 : id 41 "Se rompe el estilo..."
# ---------------------------------------------------------------------------
$ws.Cells.Item(41,2).Copy() | Out-Null
$ws.Cells.Item(43,2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(43,2).Value = "Se rompe el estilo cuando se selecciona la fecha (FECHA ELIMINATORIA)"

$ws.Cells.Item(43,4).Value = "Flor"
$ws.Cells.Item(43,5).Value = "torneo/fechas"

$ws.Cells.Item(40,6).Copy() | Out-Null
$ws.Cells.Item(43,6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(43,6).Value = "PENDIENTE"

$ws.Rows.Item(43).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 44 : id 42 "Filtro de equipos en fixture"
# ---------------------------------------------------------------------------
$ws.Cells.Item(41,1).Copy() | Out-Null
$ws.Cells.Item(44,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(44,1).Value = 42

$ws.Cells.Item(41,2).Copy() | Out-Null
$ws.Cells.Item(44,2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(44,2).Value = "Filtro de equipos en fixture"

$ws.Cells.Item(41,3).Copy() | Out-Null
$ws.Cells.Item(44,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(44,3).Value = "Que aplique para la fase correcta. Y si es fase eliminatoria Que no aparezca"

$ws.Cells.Item(44,4).Value = "Tony"
$ws.Cells.Item(44,5).Value = "torneo/fixture"

$ws.Cells.Item(41,6).Copy() | Out-Null
$ws.Cells.Item(44,6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(44,6).Value = "PENDIENTE"

$ws.Rows.Item(44).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 45 : id 43 (empty placeholder row, like rows 42/43 were before the edit)
# ---------------------------------------------------------------------------
$ws.Cells.Item(41,1).Copy() | Out-Null
$ws.Cells.Item(45,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(45,1).Value = 43

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# View state: scroll and active selection, like the edited workbook.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("C43").Select()
